# "code clean and generic" -- the Scheduler.xlsx "Sceduler-Detail" sheet
# listed "End Date" and "Recurrence" as required ("YES"); flip those two
# to "NO" so the sheet no longer hard-codes those fields as mandatory.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sceduler-Detail")

# End Date -> Required = NO
$ws.Range("B5").Value = "NO"

# Recurrence -> Required = NO
$ws.Range("B8").Value = "NO"

# Leave the cursor on B7, matching the saved selection state.
$ws.Range("B7").Select() | Out-Null
